# Add a new "double" column (F) with numeric test data, matching the
# existing header styling/formatting used by the other header cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell F1: same (bold) style as the other header cells.
$ws.Range("F1").Value = "double"
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats

# New data values in column F.
$ws.Range("F2").Value = 1.05
$ws.Range("F3").Value = 2.5
$ws.Range("F4").Value = 3.14

# Give column F an explicit width, close to the other custom column (D).
$ws.Columns("F").ColumnWidth = 13

# Reflect the last-entered cell as the active selection.
$ws.Range("F2").Select()
